$d = $word.ActiveDocument

# --- Helper text (curly apostrophe + common paragraph body) ---
$apos = [char]0x2019
$bodyText = "To take advantage of this template" + $apos + "s design, use the Styles gallery on the Home tab. You can format your headings by using heading styles, or highlight important text using other styles, like Emphasis and Intense Quote. These styles come in formatted to look great and work together to help communicate your ideas."

# --- 1. Merge the two runs of the second "To take advantage..." paragraph
#        (originally split "...format yo" | "ur headings...") into one run ---
$d.Content.Find.Execute("format your headings", $true, $false, $false, $false, $false, $true, 1, $false, "format your headings", 2) | Out-Null

# --- 2. Merge the two runs of the first "To take advantage..." paragraph
#        (originally split "...Emphasis and " | "Intense Quote...") into one run ---
$d.Content.Find.Execute("Emphasis and Intense Quote", $true, $false, $false, $false, $false, $true, 1, $false, "Emphasis and Intense Quote", 2) | Out-Null

# --- 3. Insert the new "Go ahead.../empty/Text Signatures/To take advantage.../Go ahead..." block
#        right after the first "Go ahead and get started." paragraph (currently paragraph 4) ---
$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.Collapse(0)
$r.InsertAfter("`rText Signatures`r" + $bodyText + "`rGo ahead and get started.`r")

# Give the newly-inserted "Text Signatures" paragraph the Heading1 style.
$d.Paragraphs.Item(6).Style = "Heading1"

# --- 4. Fix up the second "Text Signatures" heading (now paragraph 10): drop the stale
#        lastRenderedPageBreak and split "Text Signatures" into "Text" / " Signatures"
#        around a (re-inserted) _GoBack bookmark. ---
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Text = "Text Signatures"

$start = $p10.Range.Start
$bm = $d.Range($start + 4, $start + 4)
$d.Bookmarks.Add("_GoBack", $bm)
